$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 0.8419516666666667
$ws.Range("H2").Value = 2.525855
$ws.Range("I2").Value = 0.005772335854373203
$ws.Range("J2").Value = 0.005772335854373203
$ws.Range("M2").Value = 116.7602843333333
$ws.Range("N2").Value = 350.280853
$ws.Range("O2").Value = 0.2375147163683267
$ws.Range("P2").Value = 0.2375147163683267
$ws.Range("Q2").Value = 98.30651599492388
$ws.Range("R2").Value = 884.7586439543149
$ws.Range("S2").Value = 0.001371014713234174
$ws.Range("T2").Value = 0.001371014713234174
$ws.Range("G3").Value = 0.8419516666666667
$ws.Range("H3").Value = 2.525855
$ws.Range("I3").Value = 0.005772335854373203
$ws.Range("J3").Value = 0.005772335854373203
$ws.Range("N3").Value = 687.4430540000001
$ws.Range("O3").Value = 0.4661340766752853
$ws.Range("P3").Value = 0.4661340766752853
$ws.Range("Q3").Value = 192.9312750179078
$ws.Range("R3").Value = 1736.38147516117
$ws.Range("S3").Value = 0.002690682443737897
$ws.Range("T3").Value = 0.002690682443737897
$ws.Range("G4").Value = 0.8419516666666667
$ws.Range("H4").Value = 2.525855
$ws.Range("I4").Value = 0.005772335854373203
$ws.Range("J4").Value = 0.005772335854373203
$ws.Range("O4").Value = 0.2963512069563879
$ws.Range("P4").Value = 0.2963512069563879
$ws.Range("Q4").Value = 122.6587350553667
$ws.Range("R4").Value = 1103.9286154983
$ws.Range("S4").Value = 0.001710638697401131
$ws.Range("T4").Value = 0.001710638697401132
$ws.Range("I5").Value = 0.8411037170617888
$ws.Range("J5").Value = 0.8411037170617888
$ws.Range("M5").Value = 116.7602843333333
$ws.Range("N5").Value = 350.280853
$ws.Range("O5").Value = 0.2375147163683267
$ws.Range("P5").Value = 0.2375147163683267
$ws.Range("Q5").Value = 14324.52617116528
$ws.Range("R5").Value = 128920.7355404876
$ws.Range("S5").Value = 0.1997745107942761
$ws.Range("T5").Value = 0.1997745107942761
$ws.Range("I6").Value = 0.8411037170617888
$ws.Range("J6").Value = 0.8411037170617888
$ws.Range("N6").Value = 687.4430540000001
$ws.Range("O6").Value = 0.4661340766752853
$ws.Range("P6").Value = 0.4661340766752853
$ws.Range("Q6").Value = 28112.57290791396
$ws.Range("S6").Value = 0.3920671045407473
$ws.Range("T6").Value = 0.3920671045407473
$ws.Range("I7").Value = 0.8411037170617888
$ws.Range("J7").Value = 0.8411037170617888
$ws.Range("O7").Value = 0.2963512069563879
$ws.Range("P7").Value = 0.2963512069563879
$ws.Range("S7").Value = 0.2492621017267653
$ws.Range("T7").Value = 0.2492621017267653
$ws.Range("I8").Value = 0.1531239470838381
$ws.Range("J8").Value = 0.1531239470838381
$ws.Range("M8").Value = 116.7602843333333
$ws.Range("N8").Value = 350.280853
$ws.Range("O8").Value = 0.2375147163683267
$ws.Range("P8").Value = 0.2375147163683267
$ws.Range("Q8").Value = 2607.797282238659
$ws.Range("R8").Value = 23470.17554014793
$ws.Range("S8").Value = 0.03636919086081647
$ws.Range("T8").Value = 0.03636919086081647
$ws.Range("I9").Value = 0.1531239470838381
$ws.Range("J9").Value = 0.1531239470838381
$ws.Range("N9").Value = 687.4430540000001
$ws.Range("O9").Value = 0.4661340766752853
$ws.Range("P9").Value = 0.4661340766752853
$ws.Range("Q9").Value = 5117.927835795936
$ws.Range("R9").Value = 46061.35052216342
$ws.Range("S9").Value = 0.07137628969080011
$ws.Range("T9").Value = 0.07137628969080011
$ws.Range("I10").Value = 0.1531239470838381
$ws.Range("J10").Value = 0.1531239470838381
$ws.Range("O10").Value = 0.2963512069563879
$ws.Range("P10").Value = 0.2963512069563879
$ws.Range("S10").Value = 0.04537846653222149
$ws.Range("T10").Value = 0.0453784665322215
